$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 7267
$ws.Range("F7").Value = 4713
$ws.Range("F8").Value = 6949
$ws.Range("F10").Value = 260
$ws.Range("F11").Value = 1468
$ws.Range("F12").Value = 847
$ws.Range("F13").Value = 144
$ws.Range("F17").Value = 148
$ws.Range("F19").Value = 212
$ws.Range("F21").Value = 1127
$ws.Range("F23").Value = 44
$ws.Range("F24").Value = 1208
$ws.Range("F29").Value = 154
$ws.Range("F32").Value = 73
$ws.Range("F35").Value = 541
$ws.Range("F37").Value = 65
$ws.Range("F38").Value = 61
$ws.Range("F39").Value = 353
$ws.Range("F41").Value = 567

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 31
$ws.Range("F14").Value = 10
$ws.Range("F26").Value = 628
$ws.Range("F28").Value = 21
$ws.Range("F31").Value = 844
$ws.Range("F32").Value = 980
$ws.Range("F33").Value = 601
$ws.Range("F39").Value = 104
$ws.Range("F40").Value = 138
$ws.Range("F42").Value = 8

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 852
$ws.Range("F6").Value = 654
$ws.Range("F8").Value = 1506
$ws.Range("F9").Value = 2376

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 852
$ws.Range("F7").Value = 654
$ws.Range("F8").Value = 654
$ws.Range("F10").Value = 7267
$ws.Range("F12").Value = 4713
$ws.Range("F13").Value = 31
$ws.Range("F14").Value = 6949
$ws.Range("F15").Value = 260
$ws.Range("F16").Value = 1468
$ws.Range("F17").Value = 10
$ws.Range("F19").Value = 847
$ws.Range("F22").Value = 148
$ws.Range("F23").Value = 1127
$ws.Range("F24").Value = 628
$ws.Range("F26").Value = 44
$ws.Range("F28").Value = 21
$ws.Range("F30").Value = 844
$ws.Range("F32").Value = 73
$ws.Range("F34").Value = 980
$ws.Range("F35").Value = 541
$ws.Range("F36").Value = 601
$ws.Range("F38").Value = 65
$ws.Range("F39").Value = 61
$ws.Range("F41").Value = 353
$ws.Range("F42").Value = 567
$ws.Range("F44").Value = 104
